$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value2 = 1635.5714
$ws.Range("I33").Value2 = 1866.5
$ws.Range("J33").Value2 = 250
$ws.Range("K33").Value2 = 1866.5
$ws.Range("L33").Value2 = 250
$ws.Range("M33").Value2 = -1637.5
$ws.Range("N33").Value2 = -708
# Row 41
$ws.Range("H41").Value2 = 780.1429000000001
$ws.Range("I41").Value2 = 466.66666
$ws.Range("J41").Value2 = 1015.25
$ws.Range("K41").Value2 = 466.66666
$ws.Range("L41").Value2 = 1015.25
$ws.Range("M41").Value2 = -26.66665999999998
$ws.Range("N41").Value2 = -1895.25
# Row 111
$ws.Range("H111").Value2 = 9259
$ws.Range("I111").Value2 = 4649.8887
$ws.Range("J111").Value2 = 30000
$ws.Range("K111").Value2 = 13949.6661
$ws.Range("L111").Value2 = 90000
$ws.Range("M111").Value2 = -10882.6661
$ws.Range("N111").Value2 = -96134
# Row 120
$ws.Range("H120").Value2 = 29000
$ws.Range("J120").Value2 = 29000
$ws.Range("L120").Value2 = 29000
$ws.Range("N120").Value2 = -38676
# Row 141
$ws.Range("H141").Value2 = 9142.380999999999
$ws.Range("I141").Value2 = 4175.294
$ws.Range("J141").Value2 = 30252.5
$ws.Range("K141").Value2 = 12525.882
$ws.Range("L141").Value2 = 90757.5
$ws.Range("M141").Value2 = -7345.882
$ws.Range("N141").Value2 = -101117.5

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value2 = 5559.1177
$ws.Range("I132").Value2 = 2352.8235
$ws.Range("K132").Value2 = 7058.470499999999
$ws.Range("M132").Value2 = -4528.470499999999

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value2 = 886
$ws.Range("I5").Value2 = 621.6
$ws.Range("J5").Value2 = 1326.6666
$ws.Range("K5").Value2 = 621.6
$ws.Range("L5").Value2 = 1326.6666
$ws.Range("M5").Value2 = -508.6
$ws.Range("N5").Value2 = -1552.6666
# Row 22
$ws.Range("H22").Value2 = 1326.4667
$ws.Range("I22").Value2 = 1408.0834
$ws.Range("J22").Value2 = 1000
$ws.Range("K22").Value2 = 1408.0834
$ws.Range("L22").Value2 = 1000
$ws.Range("M22").Value2 = -1235.0834
$ws.Range("N22").Value2 = -1346
# Row 86
$ws.Range("H86").Value2 = 1900
$ws.Range("I86").Value2 = 1850
$ws.Range("J86").Value2 = 2000
$ws.Range("K86").Value2 = 1850
$ws.Range("L86").Value2 = 2000
$ws.Range("M86").Value2 = -727
$ws.Range("N86").Value2 = -4246
# Row 89
$ws.Range("H89").Value2 = 1900
$ws.Range("I89").Value2 = 1850
$ws.Range("J89").Value2 = 2000
$ws.Range("K89").Value2 = 9250
$ws.Range("L89").Value2 = 10000
$ws.Range("M89").Value2 = -3634
$ws.Range("N89").Value2 = -21232
# Row 134
$ws.Range("H134").Value2 = 6311.1113
$ws.Range("J134").Value2 = 8086.091
$ws.Range("L134").Value2 = 24258.273
$ws.Range("N134").Value2 = -29328.273

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 1808.3091
$ws.Range("I31").Value2 = 1353.2709
$ws.Range("J31").Value2 = 4928.5713
$ws.Range("K31").Value2 = 1353.2709
$ws.Range("L31").Value2 = 4928.5713
$ws.Range("M31").Value2 = -1058.2709
$ws.Range("N31").Value2 = -5518.5713
# Row 34
$ws.Range("H34").Value2 = 1808.3091
$ws.Range("I34").Value2 = 1353.2709
$ws.Range("J34").Value2 = 4928.5713
$ws.Range("K34").Value2 = 1353.2709
$ws.Range("L34").Value2 = 4928.5713
$ws.Range("M34").Value2 = -1151.2709
$ws.Range("N34").Value2 = -5332.5713
# Row 132
$ws.Range("H132").Value2 = 2311.0334
$ws.Range("I132").Value2 = 1250.6428
$ws.Range("J132").Value2 = 3238.875
$ws.Range("K132").Value2 = 3751.9284
$ws.Range("L132").Value2 = 9716.625
$ws.Range("M132").Value2 = -1221.9284
$ws.Range("N132").Value2 = -14776.625

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 840.9524
$ws.Range("I5").Value2 = 310.58823
$ws.Range("K5").Value2 = 931.76469
$ws.Range("M5").Value2 = -819.76469
# Row 122
$ws.Range("H122").Value2 = 3469.9587
$ws.Range("I122").Value2 = 373.85715
$ws.Range("J122").Value2 = 3710.7666
$ws.Range("K122").Value2 = 3364.71435
$ws.Range("L122").Value2 = 33396.8994
$ws.Range("M122").Value2 = -914.7143499999997
$ws.Range("N122").Value2 = -38296.8994
# Row 135
$ws.Range("H135").Value2 = 840.9524
$ws.Range("I135").Value2 = 310.58823
$ws.Range("K135").Value2 = 2795.29407
$ws.Range("M135").Value2 = -260.2940699999999
# Row 140
$ws.Range("H140").Value2 = 1096.8462
$ws.Range("I140").Value2 = 1096.8462
$ws.Range("K140").Value2 = 3290.5386
$ws.Range("M140").Value2 = 1889.4614

$ws = $wb.Worksheets.Item("GSM")
# Row 120
$ws.Range("H120").Value2 = 23333.334
$ws.Range("J120").Value2 = 23333.334
$ws.Range("L120").Value2 = 23333.334
$ws.Range("N120").Value2 = -33009.334
# Row 132
$ws.Range("H132").Value2 = 2453557.5
$ws.Range("I132").Value2 = 4632053
$ws.Range("K132").Value2 = 13896159
$ws.Range("M132").Value2 = -13893629

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value2 = 866.92
$ws.Range("I22").Value2 = 661.2727
$ws.Range("J22").Value2 = 1028.5
$ws.Range("K22").Value2 = 661.2727
$ws.Range("L22").Value2 = 1028.5
$ws.Range("M22").Value2 = -366.2727
$ws.Range("N22").Value2 = -1618.5
# Row 27
$ws.Range("H27").Value2 = 866.92
$ws.Range("I27").Value2 = 661.2727
$ws.Range("J27").Value2 = 1028.5
$ws.Range("K27").Value2 = 661.2727
$ws.Range("L27").Value2 = 1028.5
$ws.Range("M27").Value2 = -554.2727
$ws.Range("N27").Value2 = -1242.5
# Row 46
$ws.Range("H46").Value2 = 477196.16
$ws.Range("I46").Value2 = 772.2222
$ws.Range("J46").Value2 = 834514.0600000001
$ws.Range("K46").Value2 = 772.2222
$ws.Range("L46").Value2 = 834514.0600000001
$ws.Range("M46").Value2 = -584.2222
$ws.Range("N46").Value2 = -834890.0600000001
# Row 121
$ws.Range("H121").Value2 = 39991
$ws.Range("J121").Value2 = 39991
$ws.Range("L121").Value2 = 39991
$ws.Range("N121").Value2 = -43485
# Row 132
$ws.Range("H132").Value2 = 9175.25
$ws.Range("I132").Value2 = 10234
$ws.Range("J132").Value2 = 5999
$ws.Range("K132").Value2 = 30702
$ws.Range("L132").Value2 = 17997
$ws.Range("M132").Value2 = -28172
$ws.Range("N132").Value2 = -23057
# Row 136
$ws.Range("H136").Value2 = 1388.2222
$ws.Range("I136").Value2 = 1499.75
$ws.Range("K136").Value2 = 4499.25
$ws.Range("M136").Value2 = -1949.25

$ws = $wb.Worksheets.Item("WVR")
# Row 121
$ws.Range("H121").Value2 = 0
$ws.Range("J121").Value2 = 0
$ws.Range("L121").Value2 = 0
$ws.Range("N121").ClearContents()
# Row 132
$ws.Range("H132").Value2 = 2552.5278
$ws.Range("I132").Value2 = 2154
$ws.Range("K132").Value2 = 6462
$ws.Range("M132").Value2 = -3932

